$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a bare number/decimal string that Excel would
# otherwise auto-convert to a numeric type; force them to keep Text format
# so the value round-trips as a string, matching the source data.
$textForceCells = @('D5', 'D6', 'D7', 'D9', 'D10', 'D12', 'D14', 'D15', 'D17', 'D18', 'D21', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D32', 'D33', 'D34', 'D35', 'D36', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D46', 'D51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '51.557.35'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '3.025.13'
$ws.Range('E3').Value = '  +3.42%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '379.13'
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('D6').Value = '102.65'
$ws.Range('E6').Value = '  +3.14%  '
$ws.Range('D7').Value = '0.544'
$ws.Range('E7').Value = '  +1.65%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.593'
$ws.Range('E9').Value = '  +3.76%  '
$ws.Range('D10').Value = '36.69'
$ws.Range('E10').Value = '  +2.91%  '
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('D12').Value = '0.0857'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = '3.497.57'
$ws.Range('E13').Value = '  +3.26%  '
$ws.Range('D14').Value = '18.49'
$ws.Range('E14').Value = '  +2.55%  '
$ws.Range('D15').Value = '7.75'
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('D16').Value = '3.014.81'
$ws.Range('E16').Value = '  +3.13%  '
$ws.Range('D17').Value = '0.984'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = '10.31'
$ws.Range('E18').Value = '  -14.30%  '
$ws.Range('D19').Value = '51.577.40'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('E20').Value = '  +1.73%  '
$ws.Range('D21').Value = '12.45'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('D22').Value = '0.0₃0961'
$ws.Range('E22').Value = '  +1.94%  '
$ws.Range('D23').Value = '70.06'
$ws.Range('E23').Value = '  +0.99%  '
$ws.Range('D24').Value = '267.84'
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('D25').Value = '3.15'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('D26').Value = '8.17'
$ws.Range('E26').Value = '  +4.14%  '
$ws.Range('D27').Value = '7.45'
$ws.Range('E27').Value = '  +5.82%  '
$ws.Range('E28').Value = '  +6.92%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '26.33'
$ws.Range('E29').Value = '  +3.74%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('D32').Value = '10.30'
$ws.Range('E32').Value = '  +3.45%  '
$ws.Range('D33').Value = '34.19'
$ws.Range('E33').Value = '  +3.31%  '
$ws.Range('D34').Value = '50.52'
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('D35').Value = '2.06'
$ws.Range('E35').Value = '  +0.58%  '
$ws.Range('D36').Value = '0.0452'
$ws.Range('E36').Value = '  +5.42%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  +6.97%  '
$ws.Range('D39').Value = '17.37'
$ws.Range('E39').Value = '  +6.49%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '1.87'
$ws.Range('E40').Value = '  +4.68%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = '0.283'
$ws.Range('E41').Value = '  +10.29%  '
$ws.Range('D42').Value = '2.58'
$ws.Range('E42').Value = '  +7.30%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').Value = '127.28'
$ws.Range('E43').Value = '  +3.23%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '0.116'
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('E45').Value = '  +9.83%  '
$ws.Range('D46').Value = '22.01'
$ws.Range('E46').Value = '  +5.61%  '
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('D49').Value = '2.028.74'
$ws.Range('E49').Value = '  +1.57%  '
$ws.Range('D50').Value = '3.322.93'
$ws.Range('E50').Value = '  +3.38%  '
$ws.Range('D51').Value = '0.0321'
$ws.Range('E51').Value = '  +2.27%  '

Write-Output "done"
